# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF holds the game "Date" as literal text (e.g. "6-8-2013-14");
# correct it to the properly formatted ISO-ish text "2014-06-08".
#
# A leading apostrophe forces Excel to keep the input as literal text
# (quote-prefixed) instead of auto-converting the yyyy-mm-dd-looking
# string into a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "'2014-06-08"
}
